$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.911.25'
$ws.Range("E2").Value = '  +1.42%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.710.59'
$ws.Range("E3").Value = '  +1.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.07'
$ws.Range("E5").Value = '  +2.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9967'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3761'
$ws.Range("E7").Value = '  +1.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.79'
$ws.Range("E8").Value = '  +2.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3470'
$ws.Range("E9").Value = '  +0.45%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.217'
$ws.Range("E10").Value = '  +2.73%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07588'
$ws.Range("E11").Value = '  +4.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9991'
$ws.Range("E12").Value = '  +0.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.38'
$ws.Range("E13").Value = '  +4.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.366'
$ws.Range("E14").Value = '  +3.32%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.101'
$ws.Range("E15").Value = '  +4.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.709.48'
$ws.Range("E16").Value = '  +1.88%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001137'
$ws.Range("E17").Value = '  +1.86%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06733'
$ws.Range("E18").Value = '  +0.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9974'
$ws.Range("E19").Value = '  +0.12%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '85.53'
$ws.Range("E20").Value = '  +4.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.44'
$ws.Range("E21").Value = '  +5.53%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.426'
$ws.Range("E22").Value = '  +4.96%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.35'
$ws.Range("E23").Value = '  +10.91%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.924.87'
$ws.Range("E24").Value = '  +1.78%  '

# Row 25
$ws.Range("E25").Value = '  +0.56%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.822'
$ws.Range("E26").Value = '  +5.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.61'
$ws.Range("E27").Value = '  +4.82%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.84'
$ws.Range("E28").Value = '  -0.63%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '133.87'
$ws.Range("E29").Value = '  +5.12%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.259'
$ws.Range("E30").Value = '  +28.36%  '

# Row 31
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.900.02'
$ws.Range("E31").Value = '  +2.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.981'
$ws.Range("E32").Value = '  +9.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.238'
$ws.Range("E33").Value = '  +5.13%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.03'
$ws.Range("E34").Value = '  +11.70%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.819'
$ws.Range("E35").Value = '  +7.31%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08868'
$ws.Range("E36").Value = '  +4.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.678'
$ws.Range("E37").Value = '  +5.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.389'
$ws.Range("E38").Value = '  +3.98%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06704'
$ws.Range("E39").Value = '  +2.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02427'
$ws.Range("E40").Value = '  +3.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2251'
$ws.Range("E41").Value = '  +6.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.289'
$ws.Range("E42").Value = '  +1.26%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6497'
$ws.Range("E43").Value = '  +4.58%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9974'
$ws.Range("E44").Value = '  +0.11%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.99'
$ws.Range("E45").Value = '  +6.58%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6204'
$ws.Range("E46").Value = '  +3.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.842'
$ws.Range("E47").Value = '  +1.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.148'
$ws.Range("E48").Value = '  +5.30%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.77'
$ws.Range("E49").Value = '  +2.79%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07330'
$ws.Range("E50").Value = '  +0.88%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.57'
$ws.Range("E51").Value = '  +6.00%  '
